# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.153.92'
$ws.Range('E2').Value = '  -0.30%  '

$ws.Range('D3').Value = '1.584.31'
$ws.Range('E3').Value = '  +0.04%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.23'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  +0.39%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.245'
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0605'
$ws.Range('E9').Value = '  -0.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.20'
$ws.Range('E10').Value = '  -1.88%  '

$ws.Range('E11').Value = '  +0.43%  '

$ws.Range('D12').Value = '1.809.62'
$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('D13').Value = '1.586.24'
$ws.Range('E13').Value = '  +1.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.01'
$ws.Range('E14').Value = '  -1.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.517'
$ws.Range('E15').Value = '  +0.12%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.04'
$ws.Range('E16').Value = '  -0.62%  '

$ws.Range('D17').Value = '26.181.33'
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('D18').Value = '0.0₃0724'
$ws.Range('E18').Value = '  -0.54%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.33'
$ws.Range('E19').Value = '  -0.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '212.40'
$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('E21').Value = '  +0.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.23'
$ws.Range('E22').Value = '  -0.46%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.16'
$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.96'
$ws.Range('E24').Value = '  +1.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.40'
$ws.Range('E25').Value = '  -0.65%  '

$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.97'
$ws.Range('E27').Value = '  -0.65%  '

$ws.Range('E28').Value = '  -0.62%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.12'
$ws.Range('E29').Value = '  -1.10%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0495'
$ws.Range('E30').Value = '  -2.14%  '

$ws.Range('E31').Value = '  +0.56%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.19'
$ws.Range('E32').Value = '  -1.01%  '

$ws.Range('D33').Value = '1.340.27'
$ws.Range('E33').Value = '  +4.32%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -1.96%  '

$ws.Range('E35').Value = '  +0.11%  '

$ws.Range('E36').Value = '  -1.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.579'
$ws.Range('E37').Value = '  -3.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0166'
$ws.Range('E38').Value = '  -0.28%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.815'
$ws.Range('E39').Value = '  +0.36%  '

$ws.Range('E40').Value = '  +3.40%  '

$ws.Range('E41').Value = '  +0.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.957'
$ws.Range('E42').Value = '  -15.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.767'
$ws.Range('E43').Value = '  +0.64%  '

$ws.Range('E44').Value = '  +0.00%  '

$ws.Range('D45').Value = '1.721.71'
$ws.Range('E45').Value = '  +0.35%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.93'
$ws.Range('E46').Value = '  -2.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.91'
$ws.Range('E47').Value = '  -3.16%  '

$ws.Range('E48').Value = '  +4.45%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.47'
$ws.Range('E49').Value = '  -2.10%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0982'
$ws.Range('E50').Value = '  -2.20%  '

$ws.Range('E51').Value = '  -1.12%  '
